# Insert a new "derivatives_dir" column between the existing "metadata_dir"
# column (F) and the "modality0" column (old G, now shifts to H).
#
# This mirrors what a user did in Excel: select column G, insert a blank
# column (pushing modality0 / modality0.input_source / modality1 /
# modality1.input_source from G:J to H:K), then fill in the new column's
# header and its one data value, and finally leave the selection on G2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("G:G").Insert() | Out-Null

$ws.Range("G1").Value = "derivatives_dir"
$ws.Range("G2").Value = "/allen/aind/stage/fake/derivatives_dir"

$ws.Range("G2").Select() | Out-Null
